$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.558.81'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '1.636.16'
$ws.Range('E3').Value = '  +0.30%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''213.75'
$ws.Range('E5').Value = '  +0.41%  '
$ws.Range('E6').Value = '  +1.96%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '''0.251'
$ws.Range('E8').Value = '  -0.56%  '
$ws.Range('E9').Value = '  +0.23%  '
$ws.Range('D10').Value = '''18.91'
$ws.Range('E10').Value = '  -0.33%  '
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').Value = '1.862.27'
$ws.Range('E12').Value = '  +0.20%  '
$ws.Range('D13').Value = '1.653.66'
$ws.Range('E13').Value = '  +1.61%  '
$ws.Range('D14').Value = '''4.15'
$ws.Range('E14').Value = '  +1.69%  '
$ws.Range('D15').Value = '''0.526'
$ws.Range('E15').Value = '  +0.06%  '
$ws.Range('D16').Value = '''65.23'
$ws.Range('E16').Value = '  +3.51%  '
$ws.Range('D17').Value = '26.573.48'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').Value = '0.0₃0743'
$ws.Range('E18').Value = '  +0.68%  '
$ws.Range('D19').Value = '''215.88'
$ws.Range('E19').Value = '  +3.15%  '
$ws.Range('E21').Value = '  +0.70%  '
$ws.Range('D22').Value = '''6.28'
$ws.Range('E22').Value = '  +1.35%  '
$ws.Range('E23').Value = '  +16.60%  '
$ws.Range('D24').Value = '''9.35'
$ws.Range('E24').Value = '  -0.90%  '
$ws.Range('D25').Value = '''146.85'
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('E27').Value = '  -0.92%  '
$ws.Range('E28').Value = '  +0.78%  '
$ws.Range('D29').Value = '''15.67'
$ws.Range('E29').Value = '  +1.97%  '
$ws.Range('D30').Value = '''0.0514'
$ws.Range('E30').Value = '  -1.49%  '
$ws.Range('E31').Value = '  -0.69%  '
$ws.Range('E32').Value = '  +3.30%  '
$ws.Range('D33').Value = '''2.98'
$ws.Range('E33').Value = '  +1.26%  '
$ws.Range('D34').Value = '1.268.54'
$ws.Range('E34').Value = '  +8.82%  '
$ws.Range('E35').Value = '  +0.87%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  +3.44%  '
$ws.Range('D38').Value = '''0.511'
$ws.Range('E38').Value = '  +1.55%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = '''0.800'
$ws.Range('E39').Value = '  -0.91%  '
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D40').Value = '''1.00'
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('E41').Value = '  -1.96%  '
$ws.Range('D42').Value = '''0.798'
$ws.Range('E42').Value = '  +0.64%  '
$ws.Range('D43').Value = '''5.36'
$ws.Range('E43').Value = '  -0.37%  '
$ws.Range('D44').Value = '1.772.50'
$ws.Range('E44').Value = '  +0.23%  '
$ws.Range('D45').Value = '''93.48'
$ws.Range('E45').Value = '  +1.47%  '
$ws.Range('E46').Value = '  +3.34%  '
$ws.Range('D47').Value = '''55.09'
$ws.Range('E47').Value = '  +0.82%  '
$ws.Range('E48').Value = '  -1.74%  '
$ws.Range('E49').Value = '  +0.40%  '
$ws.Range('D50').Value = '''7.59'
$ws.Range('E50').Value = '  +0.32%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '''0.0962'
$ws.Range('E51').Value = '  +2.54%  '
